# Updated cryptos list — refreshes Price (column D) and Volume(1h) (column E)
# values for each coin row, plus a few coin-identity swaps (rows 33/34,
# 39/40, 42/43 exchanged their Coin/Link/Price/Volume content).
#
# Some Price values look numeric (e.g. "0.998", "306.34") but must remain
# plain text to match the original sheet's inlineStr cells, so we force the
# cell's number format to Text ("@") before assigning those. Values that are
# not parseable as plain numbers (e.g. "42.807.63", "0.0₃0990") are safe
# to assign directly -- Excel keeps them as text already.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.807.63'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.539.57'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.34'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.63'
$ws.Range('E6').Value = '  +8.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.585'
$ws.Range('E7').Value = '  +1.89%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.550'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.37'
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0819'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.68'
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').Value = '2.924.82'
$ws.Range('E14').Value = '  -0.81%  '
$ws.Range('D15').Value = '2.564.56'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.21'
$ws.Range('E16').Value = '  +6.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.868'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '42.913.14'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.09'
$ws.Range('E19').Value = '  +3.08%  '
$ws.Range('D20').Value = '0.0₃0990'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.52'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.74'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '254.46'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.94'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.40'
$ws.Range('E26').Value = '  -4.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.54'
$ws.Range('E28').Value = '  +2.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.33'
$ws.Range('E29').Value = '  +9.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.94'
$ws.Range('E30').Value = '  +5.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.19'
$ws.Range('E31').Value = '  +2.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.73'
$ws.Range('E32').Value = '  +3.06%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.34'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.12'
$ws.Range('E34').Value = '  -1.24%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.42'
$ws.Range('E37').Value = '  +2.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.116'
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.120'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.29'
$ws.Range('E40').Value = '  +3.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.48'
$ws.Range('E41').Value = '  +2.19%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.93'
$ws.Range('E42').Value = '  +1.27%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.09'
$ws.Range('E43').Value = '  +0.49%  '
$ws.Range('E44').Value = '  -1.85%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '2.044.63'
$ws.Range('E46').Value = '  -1.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.22'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.05'
$ws.Range('E48').Value = '  -2.40%  '
$ws.Range('D49').Value = '2.781.89'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.84'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('E51').Value = '  +0.93%  '
